$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.653.15'
$ws.Range('E2').Value = '  +4.04%  '
$ws.Range('D3').Value = '3.409.05'
$ws.Range('E3').Value = '  +2.48%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.55%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '595.93'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +7.41%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '188.91'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.48%  '
$ws.Range('E7').Value = '  +4.04%  '
$ws.Range('E8').Value = '  +0.00%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.187'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +5.09%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '0.594'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +2.44%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '48.00'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +4.55%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.0000284'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +7.18%  '
$ws.Range('B13').Value = 'BitcoinCash'
$ws.Range('C13').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '645.85'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +11.47%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.946.12'
$ws.Range('E14').Value = '  +2.22%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '8.69'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  +2.53%  '
$ws.Range('D16').Value = '68.599.28'
$ws.Range('E16').Value = '  +3.90%  '
$ws.Range('D17').Value = '3.396.14'
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '18.20'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  +2.43%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.119'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +1.80%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '11.20'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.67%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '0.918'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  +2.81%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '18.06'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.14%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.14'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +2.43%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '100.61'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +1.93%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '4.07'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  +2.93%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.88'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +6.96%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '9.82'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +4.83%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '33.13'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  +8.39%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '8.81'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +4.70%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '6.97'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +4.74%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '616.40'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +7.30%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.91'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  +4.83%  '
$ws.Range('D33').Value = '4.048.97'
$ws.Range('E33').Value = '  +8.86%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '11.19'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +3.04%  '
$ws.Range('E35').Value = '  +3.79%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.10%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '56.65'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +2.64%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '2.83'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  +8.01%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.132'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +4.46%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '3.32'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +5.86%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '34.04'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = '0.0₃0714'
$ws.Range('E42').Value = '  +3.82%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '0.348'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '3.43'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  +1.63%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.0428'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  +4.35%  '
$ws.Range('E46').Value = '  +1.89%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.63'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  +4.39%  '
$ws.Range('E48').Value = '  +12.68%  '
$ws.Range('E49').Value = '  -0.08%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '129.23'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +1.95%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '7.86'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +7.64%  '
